# Update countries & provincias Spain
# Refresh the "Pais" covid-stats sheet: new timestamp, updated per-country
# counters, and a handful of countries that swapped rank (and therefore
# row) because their "Casos totales" changed relative to their neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 23 de Julio de 2020 a las 14:39"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 4102002
$ws.Range("C4").Value = 1127
$ws.Range("E4").Value = 2012301
$ws.Range("G4").Value = 15
$ws.Range("H4").Value = 146198

# India (row 6)
$ws.Range("B6").Value = 1241687
$ws.Range("C6").Value = 2003
$ws.Range("D6").Value = 784460
$ws.Range("E6").Value = 427321

# Bielorrusia (row 35)
$ws.Range("B35").Value = 66688
$ws.Range("C35").Value = 167
$ws.Range("D35").Value = 59439
$ws.Range("E35").Value = 6730
$ws.Range("G35").Value = 6
$ws.Range("H35").Value = 519

# Ucrania/Kuwait swap ranks (rows 38-39)
$ws.Range("A38").Value = "Kuwait"
$ws.Range("B38").Value = 61872
$ws.Range("C38").Value = 687
$ws.Range("D38").Value = 52247
$ws.Range("E38").Value = 9204
$ws.Range("G38").Value = 4
$ws.Range("H38").Value = 421

$ws.Range("A39").Value = "Ucrania"
$ws.Range("B39").Value = 61851
$ws.Range("C39").Value = 856
$ws.Range("D39").Value = 34000
$ws.Range("E39").Value = 26300
$ws.Range("G39").Value = 17
$ws.Range("H39").Value = 1551

# Paises Bajos (row 44)
$ws.Range("B44").Value = 52404
$ws.Range("C44").Value = 163

# Moldavia (row 63)
$ws.Range("D63").Value = 15174
$ws.Range("E63").Value = 5911

# El Salvador/Dinamarca swap ranks (rows 73-74)
$ws.Range("A73").Value = "Dinamarca"
$ws.Range("B73").Value = 13390
$ws.Range("C73").Value = 42
$ws.Range("D73").Value = 12299
$ws.Range("E73").Value = 479
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 612

$ws.Range("A74").Value = "El Salvador"
$ws.Range("B74").Value = 13377
$ws.Range("C74").Value = 402
$ws.Range("D74").Value = 7276
$ws.Range("E74").Value = 5729
$ws.Range("G74").Value = 9
$ws.Range("H74").Value = 372

# Republica de Africa Central/Croacia swap ranks (rows 98-99)
$ws.Range("A98").Value = "Croacia"
$ws.Range("B98").Value = 4634
$ws.Range("C98").Value = 104
$ws.Range("D98").Value = 3394
$ws.Range("E98").Value = 1112
$ws.Range("G98").Value = 3
$ws.Range("H98").Value = 128

$ws.Range("A99").Value = "Republica de Africa Central"
$ws.Range("B99").Value = 4574
$ws.Range("D99").Value = 1437
$ws.Range("E99").Value = 3080
$ws.Range("H99").Value = 57

# Sudan del Sur / Libia / Cabo Verde / Hong Kong shift ranks (rows 118-121)
$ws.Range("A118").Value = "Hong Kong"
$ws.Range("B118").Value = 2250
$ws.Range("C118").Value = 118
$ws.Range("D118").Value = 1379
$ws.Range("E118").Value = 856
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 15

$ws.Range("A119").Value = "Sudan del Sur"
$ws.Range("B119").Value = 2211
$ws.Range("D119").Value = 1175
$ws.Range("E119").Value = 991
$ws.Range("H119").Value = 45

$ws.Range("A120").Value = "Libia"
$ws.Range("B120").Value = 2176
$ws.Range("D120").Value = 489
$ws.Range("E120").Value = 1634
$ws.Range("H120").Value = 53

$ws.Range("A121").Value = "Cabo Verde"
$ws.Range("B121").Value = 2154
$ws.Range("D121").Value = 1132
$ws.Range("E121").Value = 1001
$ws.Range("H121").Value = 21

# Islandia (row 129)
$ws.Range("B129").Value = 1841
$ws.Range("C129").Value = 1
$ws.Range("D129").Value = 1823

# Vietnam (row 162)
$ws.Range("B162").Value = 412
$ws.Range("C162").Value = 4
$ws.Range("E162").Value = 47

# San Martin (Parte Holandesa) (row 190)
$ws.Range("B190").Value = 81
$ws.Range("C190").Value = 2
$ws.Range("E190").Value = 3

# Groenlandia/Islas Malvinas swap ranks (rows 210-211)
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"

